# Apply updated vm_pu results for the 380 kV case (rows 2-25, columns B-F and I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039216618482286
$ws.Range("D2").Value = 1.046357541331056
$ws.Range("E2").Value = 1.047857386747824
$ws.Range("F2").Value = 1.058917009885102
$ws.Range("I2").Value = 1.040276226565228
$ws.Range("J2").Value = 1.044309834564598
$ws.Range("K2").Value = 1.049123203024651
$ws.Range("L2").Value = 1.050618852859598
$ws.Range("M2").Value = 1.061647935920581
$ws.Range("N2").Value = 1.045792874313692

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040162936897491
$ws.Range("D3").Value = 1.047073490039572
$ws.Range("E3").Value = 1.048676432164096
$ws.Range("F3").Value = 1.059769035826908
$ws.Range("I3").Value = 1.040464516732097
$ws.Range("J3").Value = 1.04490128326365
$ws.Range("K3").Value = 1.04965102672486
$ws.Range("L3").Value = 1.051249808405838
$ws.Range("M3").Value = 1.062313993019986
$ws.Range("N3").Value = 1.046385162937737

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040775717791575
$ws.Range("D4").Value = 1.047536831201523
$ws.Range("E4").Value = 1.049207103766685
$ws.Range("F4").Value = 1.060320855789997
$ws.Range("I4").Value = 1.040584755562964
$ws.Range("J4").Value = 1.045283816347476
$ws.Range("K4").Value = 1.049991966323565
$ws.Range("L4").Value = 1.051658121554631
$ws.Range("M4").Value = 1.062744832446882
$ws.Range("N4").Value = 1.046768239262418

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041033437400736
$ws.Range("D5").Value = 1.047731636111191
$ws.Range("E5").Value = 1.049430362854265
$ws.Range("F5").Value = 1.060552959511864
$ws.Range("I5").Value = 1.040634920776059
$ws.Range("J5").Value = 1.045444590562447
$ws.Range("K5").Value = 1.05013515307142
$ws.Range("L5").Value = 1.05182978542373
$ws.Range("M5").Value = 1.062925921491314
$ws.Range("N5").Value = 1.046929241795216

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04107671587125
$ws.Range("D6").Value = 1.047764345622707
$ws.Range("E6").Value = 1.049467858654371
$ws.Range("F6").Value = 1.060591937656232
$ws.Range("I6").Value = 1.040643321243781
$ws.Range("J6").Value = 1.045471582725245
$ws.Range("K6").Value = 1.050159186262463
$ws.Range("L6").Value = 1.051858609036628
$ws.Range("M6").Value = 1.062956324996699
$ws.Range("N6").Value = 1.046956272289981

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040779161036443
$ws.Range("D7").Value = 1.047539434132218
$ws.Range("E7").Value = 1.049210086320099
$ws.Range("F7").Value = 1.060323956707275
$ws.Range("I7").Value = 1.040585427379617
$ws.Range("J7").Value = 1.045285964789535
$ws.Range("K7").Value = 1.049993880160234
$ws.Range("L7").Value = 1.051660415301722
$ws.Range("M7").Value = 1.062747252310144
$ws.Range("N7").Value = 1.046770390755512

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039536337458344
$ws.Range("D8").Value = 1.046599483202101
$ws.Range("E8").Value = 1.04813404228534
$ws.Range("F8").Value = 1.059204851083681
$ws.Range("I8").Value = 1.040340190432633
$ws.Range("J8").Value = 1.044509752746766
$ws.Range("K8").Value = 1.049301706546802
$ws.Range("L8").Value = 1.050832077670966
$ws.Range("M8").Value = 1.061873062044297
$ws.Range("N8").Value = 1.04599307640261

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037349816011797
$ws.Range("D9").Value = 1.044943807907629
$ws.Range("E9").Value = 1.046243300278029
$ws.Range("F9").Value = 1.05723676162515
$ws.Range("I9").Value = 1.039895849516591
$ws.Range("J9").Value = 1.04314068511918
$ws.Range("K9").Value = 1.048077483233626
$ws.Range("L9").Value = 1.049372821215318
$ws.Range("M9").Value = 1.060331583856112
$ws.Range("N9").Value = 1.044622064541958

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035894548019024
$ws.Range("D10").Value = 1.043840544440646
$ws.Range("E10").Value = 1.0449865150232
$ws.Range("F10").Value = 1.055927431341146
$ws.Range("I10").Value = 1.03959146198807
$ws.Range("J10").Value = 1.042227174880679
$ws.Range("K10").Value = 1.047258361595917
$ws.Range("L10").Value = 1.048400310135251
$ws.Range("M10").Value = 1.059303310466299
$ws.Range("N10").Value = 1.043707257014137

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035264986467072
$ws.Range("D11").Value = 1.043362961466952
$ws.Range("E11").Value = 1.044443211464008
$ws.Range("F11").Value = 1.055361144096552
$ws.Range("I11").Value = 1.039457732362926
$ws.Range("J11").Value = 1.041831438721865
$ws.Range("K11").Value = 1.04690298221347
$ws.Range("L11").Value = 1.04797929427613
$ws.Range("M11").Value = 1.058857926036128
$ws.Range("N11").Value = 1.043310958864587

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03503122729101
$ws.Range("D12").Value = 1.043185588027701
$ws.Range("E12").Value = 1.04424153992883
$ws.Range("F12").Value = 1.055150900744204
$ws.Range("I12").Value = 1.039407770282946
$ws.Range("J12").Value = 1.041684418654898
$ws.Range("K12").Value = 1.046770875303775
$ws.Range("L12").Value = 1.047822924598815
$ws.Range("M12").Value = 1.058692471416802
$ws.Range("N12").Value = 1.043163730012261

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035081365428607
$ws.Range("D13").Value = 1.04322363420653
$ws.Range("E13").Value = 1.044284792998052
$ws.Range("F13").Value = 1.055195994063779
$ws.Range("I13").Value = 1.039418500382921
$ws.Range("J13").Value = 1.041715956119978
$ws.Range("K13").Value = 1.046799217335553
$ws.Range("L13").Value = 1.047856465757068
$ws.Range("M13").Value = 1.058727962833574
$ws.Range("N13").Value = 1.043195312264159

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035245662059297
$ws.Range("D14").Value = 1.04334829926279
$ws.Range("E14").Value = 1.044426538457354
$ws.Range("F14").Value = 1.055343763247134
$ws.Range("I14").Value = 1.039453608372174
$ws.Range("J14").Value = 1.041819286522442
$ws.Range("K14").Value = 1.046892064313509
$ws.Range("L14").Value = 1.04796636841022
$ws.Range("M14").Value = 1.058844249885933
$ws.Range("N14").Value = 1.043298789407646

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035346902329212
$ws.Range("D15").Value = 1.043425112497197
$ws.Range("E15").Value = 1.044513890524954
$ws.Range("F15").Value = 1.055434822141145
$ws.Range("I15").Value = 1.03947520129593
$ws.Range("J15").Value = 1.04188294837446
$ws.Range("K15").Value = 1.046949256761808
$ws.Range("L15").Value = 1.048034084999347
$ws.Range("M15").Value = 1.058915895696845
$ws.Range("N15").Value = 1.043362541666797

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035936342233083
$ws.Range("D16").Value = 1.043872243072521
$ws.Range("E16").Value = 1.045022591225995
$ws.Range("F16").Value = 1.055965028030187
$ws.Range("I16").Value = 1.039600296624897
$ws.Range("J16").Value = 1.042253434878057
$ws.Range("K16").Value = 1.047281932423205
$ws.Range("L16").Value = 1.048428253521252
$ws.Range("M16").Value = 1.059332866426152
$ws.Range("N16").Value = 1.043733554303724

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036306238066284
$ws.Range("D17").Value = 1.044152754112794
$ws.Range("E17").Value = 1.045341925790383
$ws.Range("F17").Value = 1.056297790367398
$ws.Range("I17").Value = 1.039678249976035
$ws.Range("J17").Value = 1.042485783760993
$ws.Range("K17").Value = 1.047490425987815
$ws.Range("L17").Value = 1.048675529361784
$ws.Range("M17").Value = 1.059594385860171
$ws.Range("N17").Value = 1.043966233148724

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036522047853578
$ws.Range("D18").Value = 1.044316384612276
$ws.Range("E18").Value = 1.045528274277934
$ws.Range("F18").Value = 1.056491948719706
$ws.Range("I18").Value = 1.039723532794377
$ws.Range("J18").Value = 1.042621291477616
$ws.Range("K18").Value = 1.047611969690479
$ws.Range("L18").Value = 1.048819769575219
$ws.Range("M18").Value = 1.059746912582429
$ws.Range("N18").Value = 1.044101933301848

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036595642848871
$ws.Range("D19").Value = 1.044372180569702
$ws.Range("E19").Value = 1.045591828845379
$ws.Range("F19").Value = 1.056558162458849
$ws.Range("I19").Value = 1.039738941486475
$ws.Range("J19").Value = 1.042667493114836
$ws.Range("K19").Value = 1.047653401542172
$ws.Range("L19").Value = 1.048868953149861
$ws.Range("M19").Value = 1.059798917984534
$ws.Range("N19").Value = 1.044148200550693

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036266545964406
$ws.Range("D20").Value = 1.044122656561295
$ws.Range("E20").Value = 1.045307655303718
$ws.Range("F20").Value = 1.056262081516523
$ws.Range("I20").Value = 1.039669905558483
$ws.Range("J20").Value = 1.042460856722807
$ws.Range("K20").Value = 1.047468063530284
$ws.Range("L20").Value = 1.04864899812712
$ws.Range("M20").Value = 1.05956632864943
$ws.Range("N20").Value = 1.043941270711285

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035197278358734
$ws.Range("D21").Value = 1.043311587887672
$ws.Range("E21").Value = 1.044384794197548
$ws.Range("F21").Value = 1.055300246106499
$ws.Range("I21").Value = 1.03944327791509
$ws.Range("J21").Value = 1.041788859004813
$ws.Range("K21").Value = 1.046864726031983
$ws.Range("L21").Value = 1.047934004412455
$ws.Range("M21").Value = 1.05881000676168
$ws.Range("N21").Value = 1.043268318679453

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034525497110266
$ws.Range("D22").Value = 1.042801765665253
$ws.Range("E22").Value = 1.043805340042585
$ws.Range("F22").Value = 1.054696086818628
$ws.Range("I22").Value = 1.039299116726675
$ws.Range("J22").Value = 1.04136619723489
$ws.Range("K22").Value = 1.046484787367276
$ws.Range("L22").Value = 1.047484543132188
$ws.Range("M22").Value = 1.058334367826465
$ws.Range("N22").Value = 1.042845056681336

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034881572404662
$ws.Range("D23").Value = 1.043072019358688
$ws.Range("E23").Value = 1.04411244473035
$ws.Range("F23").Value = 1.055016307129939
$ws.Range("I23").Value = 1.039375697499251
$ws.Range("J23").Value = 1.041590272121755
$ws.Range("K23").Value = 1.046686256201128
$ws.Range("L23").Value = 1.047722802745751
$ws.Range("M23").Value = 1.058586522953111
$ws.Range("N23").Value = 1.043069449780242

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036284480944347
$ws.Range("D24").Value = 1.044136256308788
$ws.Range("E24").Value = 1.045323140395791
$ws.Range("F24").Value = 1.056278216612466
$ws.Range("I24").Value = 1.039673676616521
$ws.Range("J24").Value = 1.042472120232416
$ws.Range("K24").Value = 1.047478168369587
$ws.Range("L24").Value = 1.048660986424197
$ws.Range("M24").Value = 1.059579006535977
$ws.Range("N24").Value = 1.04395255021637

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037914663664939
$ws.Range("D25").Value = 1.045371754629222
$ws.Range("E25").Value = 1.046731455088358
$ws.Range("F25").Value = 1.057745085825752
$ws.Range("I25").Value = 1.040012163788802
$ws.Range("J25").Value = 1.043494767336229
$ws.Range("K25").Value = 1.048394503355794
$ws.Range("L25").Value = 1.049750021609711
$ws.Range("M25").Value = 1.060730208275392
$ws.Range("N25").Value = 1.044976649596369
